$d = $word.ActiveDocument

# Locate the trailing empty paragraph (the one right before sectPr) - this is
# where Problem 2's content will be inserted. The final paragraph of the
# inserted fragment merges into this pre-existing paragraph.
$lastIdx = $d.Paragraphs.Count
$lastP = $d.Paragraphs.Item($lastIdx)
$insertRange = $d.Range($lastP.Range.Start, $lastP.Range.Start)

$frag = @'
    <w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"/>
    <w:p>
      <w:pPr>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
      </w:pPr>
      <w:r>
        <w:rPr>
          <w:b/>
          <w:u w:val="single"/>
        </w:rPr>
        <w:t>Problem 2</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>1.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>a) You have 20 objects. These objects belong in pairs. Of these 20, you have 3 classifications. There are 5 pairs of the first class, three pairs of the second class, and two pairs of the third class. You are unable to discern between the classes while taking the objects. Keeping this in mind, what is the least number you must take in order to ensure you have:</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" w:firstLine="720"/>
      </w:pPr>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>At least one pair from the same class.</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720" w:firstLine="720"/>
      </w:pPr>
      <w:proofErr w:type="gramStart"/>
      <w:r>
        <w:t>At least one pair from each class.</w:t>
      </w:r>
      <w:proofErr w:type="gramEnd"/>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:firstLine="720"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>b) This problem is also one that has been around for some time. However, it is often told with different amounts of each classification. It is apparent that this is a simple probability problem, in which one must minimize the objective function, while adhering to the constraints.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>c) The overall goal is to come up with the minimum amount of objects you must take in order to ENSURE that you have enough to reach the goal. In other words, one must assume that the “worst-case scenario” will happen, and should pull as much as necessary to make sure that the constraints are still fulfilled.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>2.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>a) The constraints are that we</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> acquire</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> all the objects req</w:t>
      </w:r>
      <w:r>
        <w:t>uired to ensure objective completion, without foreknowledge.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> </w:t>
      </w:r>
      <w:r>
        <w:t>One must acquire this while tak</w:t>
      </w:r>
      <w:r>
        <w:t>ing the least amount of objects, while still fulfilling the first constraint.</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve"> So, the first constraint is to take a number of socks that ensures you reach the objective. The second constraint is to take the minimum amount of socks possible.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t xml:space="preserve">b) </w:t>
      </w:r>
      <w:r>
        <w:t>The sub-goals are t</w:t>
      </w:r>
      <w:r>
        <w:t>hat you take a number of socks. This is, of course, very simple, but if you were to add the constraints that the number of socks be a certain combination, then you would have the overall objective.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>3.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>a) There only seems to be one solution that matches the constraints for both objectives. However, on might attempt to take a certain number of socks that has a high likelihood of achieving the desired result, rather than one that would ensure it.</w:t>
      </w:r>
    </w:p>
    <w:p/>
    <w:p>
      <w:r>
        <w:t>4.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>a) This strategy, though possibly being likely to result in the objective, doesn’t ensure it absolutely. Therefore, it doesn’t meet the goals.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:tab/>
        <w:t>b) The only solutions that would work for all cases would be the one below.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:r>
        <w:t>5.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>a) The only solution for the first objective (acquiri</w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">ng at least one matching pair) </w:t>
      </w:r>
      <w:r>
        <w:t xml:space="preserve">would be to take four socks. As there are </w:t>
      </w:r>
      <w:r>
        <w:t>only three classes (or colors) of socks, then even if you took one of each class in the first three pulls, the fourth would be one of those three colors, matching the one already taken and fulfilling the constraint that we take the minimum amount, while ensuring we reach the objective. For the second objective, the only solution would be to take 18 socks. Though unlikely, it is possible that the first 16 socks taken could be black and brown. Therefore, you would still have to take two of the remaining to absolutely reach the objective. This is the minimum you can pull to remain certain you’ve taken the amount necessary.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
      <w:r>
        <w:t>b) As this is a simple probability problem, I did no test cases; however, one can be certain that these answers are sound, and that there are no other options.</w:t>
      </w:r>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
    </w:p>
    <w:p>
      <w:pPr>
        <w:ind w:left="720"/>
      </w:pPr>
    </w:p>
'@

$insertRange.InsertXML($frag)

# The fragment's own trailing (31st, merge) paragraph didn't carry over its
# <w:ind w:left="720"/> because InsertXML merges its final paragraph into the
# pre-existing one and keeps that paragraph's own pPr. Fix that up directly.
$finalIdx = $d.Paragraphs.Count
$finalP = $d.Paragraphs.Item($finalIdx)
$finalP.Range.ParagraphFormat.LeftIndent = 36

# Move the "_GoBack" bookmark from the end of the "b) As I had already
# known..." paragraph down onto this new final paragraph, matching the diff.
$d.Bookmarks.Add("_GoBack", $finalP.Range)
